# The edit swaps the data between row 25 and row 26 of the "Artfynd" sheet
# (two species observation records trade places), while columns C
# (Valideringsstatus) and P (Lokalnamn) keep their original values because
# they were already identical on both rows.
#
# Target values after the edit:
#   Row 25 <- former row 26 data (Knärot / Goodyera repens, VU, ...)
#   Row 26 <- former row 25 data (Kandelabersvamp / Artomyces pyxidatus, NT, ...)
#
# Additionally, the "Ålder-Stadium" value ("blomning", column K) that used to
# sit on row 26 moves to row 25 (together with a now-blank column L cell),
# while row 26's K cell becomes blank (present but empty) and its L cell
# disappears entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("Antal") holds numeric-looking text ("2", "15"); force text
# format so Excel keeps it as a string instead of converting to a number.
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I26").NumberFormat = "@"

# ---- Row 25 gets the former row 26 values ----
$ws.Range("A25").Value = 111595483
$ws.Range("B25").Value = 96348
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("I25").Value = "15"
$ws.Range("J25").Value = "plantor/tuvor"
$ws.Range("K25").Value = "blomning"
# L25 did not exist before the edit and must stay present but blank;
# touching its number format alone creates a style-only, content-less cell.
$ws.Range("L25").NumberFormat = "General"
$ws.Range("Q25").Value = 578724.2708698318
$ws.Range("R25").Value = 6410783.051849495
$ws.Range("AC25").Value = "6 blommor"

# ---- Row 26 gets the former row 25 values ----
$ws.Range("A26").Value = 111595525
$ws.Range("B26").Value = 90151
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 366
$ws.Range("F26").Value = "Kandelabersvamp"
$ws.Range("G26").Value = "Artomyces pyxidatus"
$ws.Range("H26").Value = "(Pers.) Jülich"
$ws.Range("I26").Value = "2"
$ws.Range("J26").Value = "fruktkroppar"
# K26 must stay present but become blank (not be removed); touching its
# number format keeps the cell node alive once its value is cleared.
$ws.Range("K26").NumberFormat = "General"
$ws.Range("K26").Value = ""
$ws.Range("Q26").Value = 578725.2392689644
$ws.Range("R26").Value = 6410706.376212179
$ws.Range("AC26").Value = "På gammal asplåga"

# L26 must be entirely cleared (no longer present), matching row 25's
# former, blank-L state. It has not been touched otherwise, so clearing its
# contents removes the cell node completely.
$ws.Range("L26").ClearContents()
